# "Chiffres COVID-19 Valais" daily update.
# Fills in the corrected figures for 2020-09-07..09 (rows 195-197) and the
# newly-published figures for 2020-09-10 (row 198), which was still an
# empty placeholder row in the previous upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 195 (2020-09-07): nouveaux cas positifs 9 -> 10 -------------------
$ws.Range("C195").Value = 10

# --- Row 196 (2020-09-08): nouveaux cas positifs 3 -> 8, nouvelles
#     admissions 0 -> 1, patients hospitalisés hors SI 4 -> 5 ---------------
$ws.Range("C196").Value = 8
$ws.Range("D196").Value = 1
$ws.Range("G196").Value = 5

# --- Row 197 (2020-09-09): nouveaux cas positifs 0 -> 3, nouvelles
#     sorties 0 -> 1 -----------------------------------------------------
$ws.Range("C197").Value = 3
$ws.Range("I197").Value = 1

# --- Row 198 (2020-09-10): figures now published for this day --------------
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 0
$ws.Range("E198").Value = 0
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 4
$ws.Range("I198").Value = 0
$ws.Range("L198").Value = "0"
$ws.Range("M198").Value = "0"

# Move the active cell of the frozen bottom-right pane to reflect where the
# author was working (matches the sheet's saved selection).
$ws.Range("I202").Select()
